$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------------
# This change regenerates the localization-status report for a new handoff:
#   - the old source file b4fe5b0b-4ce6-4a74-a8d0-dedd5801e610.md is replaced
#     by 2f48f729-5dd0-4636-a715-fa045eb9a4ba.md
#   - its handoff xlf filenames/hashes and handoff datetimes are refreshed
#   - the row describing df1e1583-cfe0-47a3-9e6c-9d66b7db1622.md ("Handoff
#     transform failed") is gone entirely (that file is no longer tracked)
# ----------------------------------------------------------------------------

$oldMd   = "b4fe5b0b-4ce6-4a74-a8d0-dedd5801e610.md"
$newMd   = "2f48f729-5dd0-4636-a715-fa045eb9a4ba.md"
$oldHash = "b4fe5b0b-4ce6-4a74-a8d0-dedd5801e610.529a8dda16127c2de88b6d910597100f0baa5b10"
$newHash = "2f48f729-5dd0-4636-a715-fa045eb9a4ba.3b41ba508a47a6cac67718a70e65a1174eb57093"

$removedMd = "df1e1583-cfe0-47a3-9e6c-9d66b7db1622.md"

$commitBase = "https://github.com/OpenLocalizationTest/oltest/blob/b541b7de7c7b93e5134f26a6c481ca42d415a89d"

# ============================================================================
# Sheet "Overview"
# ============================================================================
$ws = $wb.Worksheets.Item("Overview")

# Find & drop the row that references the removed file, shifting the rows
# below it up (this also fixes the sheet's used-range / dimension).
for ($r = $ws.UsedRange.Rows.Count(); $r -ge 1; $r--) {
    if ($ws.Cells.Item($r, 1).Value() -eq $removedMd) {
        $ws.Range("A" + $r).EntireRow.Delete()
    }
}

# Rename the surviving UUID reference.
for ($r = 1; $r -le $ws.UsedRange.Rows.Count(); $r++) {
    if ($ws.Cells.Item($r, 1).Value() -eq $oldMd) {
        $ws.Cells.Item($r, 1).Value = $newMd
    }
}

# Rebuild hyperlinks to match the (now-shifted) rows & renamed display text.
$ws.Hyperlinks.Delete()
for ($r = 2; $r -le $ws.UsedRange.Rows.Count(); $r++) {
    $name = $ws.Cells.Item($r, 1).Value()
    if ($name -eq $newMd) {
        $ws.Hyperlinks.Add($ws.Range("A" + $r), $commitBase + "/e2e/" + $newMd, "", "", $newMd)
    } elseif ($name -eq ".localization-config") {
        $ws.Hyperlinks.Add($ws.Range("A" + $r), $commitBase + "/.localization-config", "", "", ".localization-config")
    }
}

# ============================================================================
# Per-locale sheets ("zh-cn", "de-de")
# ============================================================================
$locales = @(
    @{ Sheet = "zh-cn"; Time = "2016-01-27 08:17:17"; OlhandoffCommit = "13908677e894c639b89ddcba417c9d73af1b2f81" },
    @{ Sheet = "de-de"; Time = "2016-01-27 08:17:29"; OlhandoffCommit = "2d714aa20eb924e1f61da6548853e64d4ad6ee55" }
)

foreach ($loc in $locales) {
    $ws = $wb.Worksheets.Item($loc.Sheet)
    $locale = $loc.Sheet
    $newXlf = $newHash + "." + $locale + ".xlf"

    # Drop the row describing the removed file (shifts later rows up).
    for ($r = $ws.UsedRange.Rows.Count(); $r -ge 1; $r--) {
        if ($ws.Cells.Item($r, 1).Value() -eq $removedMd) {
            $ws.Range("A" + $r).EntireRow.Delete()
        }
    }

    # Update the surviving data row: file name, xlf name, and handoff time.
    for ($r = 1; $r -le $ws.UsedRange.Rows.Count(); $r++) {
        if ($ws.Cells.Item($r, 1).Value() -eq $oldMd) {
            $ws.Cells.Item($r, 1).Value = $newMd
            $ws.Cells.Item($r, 3).Value = $newXlf
            $ws.Cells.Item($r, 4).Value = $loc.Time
        }
    }

    # Rebuild hyperlinks for this sheet.
    $ws.Hyperlinks.Delete()
    for ($r = 2; $r -le $ws.UsedRange.Rows.Count(); $r++) {
        $name = $ws.Cells.Item($r, 1).Value()
        if ($name -eq $newMd) {
            $ws.Hyperlinks.Add($ws.Range("A" + $r), $commitBase + "/e2e/" + $newMd, "", "", $newMd)

            $xlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/" + $loc.OlhandoffCommit + "/ol-handoff/OpenLocalizationTestOrg/oltest." + $locale + "/tianzh/" + $newXlf
            $ws.Hyperlinks.Add($ws.Range("C" + $r), $xlfUrl, "", "", $newXlf)
        } elseif ($name -eq ".localization-config") {
            $ws.Hyperlinks.Add($ws.Range("A" + $r), $commitBase + "/.localization-config", "", "", ".localization-config")
        }
    }
}
